$p = $ppt.ActivePresentation
$s = $p.Slides.Item(3)

# EMU -> points (PowerPoint COM Left/Top/Width/Height are in points)
$left   = 5280212 / 914400 * 72
$top    = 3836894 / 914400 * 72
$width  = 3702423 / 914400 * 72
$height = 2862322 / 914400 * 72

# The presentation's internal shape-id counter is monotonically increasing
# across the whole editing session (it does not reuse ids after deletes).
# The target shape must land on id=6 / name="TextBox 5", so burn two ids
# first with scratch textboxes that get removed again.
$junk1 = $s.Shapes.AddTextbox(1, 0, 0, 10, 10)
$junk1.Delete()
$junk2 = $s.Shapes.AddTextbox(1, 0, 0, 10, 10)
$junk2.Delete()

$tb = $s.Shapes.AddTextbox(1, $left, $top, $width, $height)

$tf = $tb.TextFrame
$tf.WordWrap = $true
$tr = $tf.TextRange

$title = "Problems with Std. Approaches:"
$body = "Many user agents attempt to deceive the server parsing engine in order to get specific content, i.e. pages optimized for GoogleBot, by adding specific tokens to the UA string.  A hierarchal regexp engine will be confused by such additions.  A discriminative algorithm may still classify the modified string correctly."

$tr.Text = $title + "`r" + $body

$tr.LanguageID = "en-US"

$p1 = $tr.Paragraphs(1, 1)
$p1.Font.Bold = $true
$p1.Font.Underline = $true

$p2 = $tr.Paragraphs(2, 1)
$base = $p1.Length

$run1 = "Many user agents attempt to deceive the server parsing engine in order to get specific content, i.e. pages optimized for "
$run2 = "GoogleBot"
$run3 = ", by adding specific tokens to the UA string.  A hierarchal "
$run4 = "regexp"
$run5 = " engine will be confused by such additions.  A discriminative algorithm may still classify the modified string correctly."

$pos = $base + 1
$r1 = $tr.Characters($pos, $run1.Length)
$pos += $run1.Length
$r2 = $tr.Characters($pos, $run2.Length)
$pos += $run2.Length
$r3 = $tr.Characters($pos, $run3.Length)
$pos += $run3.Length
$r4 = $tr.Characters($pos, $run4.Length)
$pos += $run4.Length
$r5 = $tr.Characters($pos, $run5.Length)

Write-Host "r1:" $r1.Text
Write-Host "r2:" $r2.Text
Write-Host "r3:" $r3.Text
Write-Host "r4:" $r4.Text
Write-Host "r5:" $r5.Text

$tb.Line.Visible = $false
